# Update the "Förändrad" date column (C) for rows 2-27 from 2023-09-03 (45172)
# to 2023-09-06 (45175), matching the automatic update recorded in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}
